# Fruta / hortaliza, semanal
# Insert two new weekly price rows for "Naranja" (Fukumoto / Thompson) at row 234,
# pushing the existing rows 234-291 down to 236-293.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows above the current row 234 (each Insert() shifts everything
# at/after row 234 down by one row).
$ws.Rows.Item(234).Insert()
$ws.Rows.Item(234).Insert()

# New row 234: Fukumoto / Primera
$ws.Cells.Item(234, 1).Value  = 11
$ws.Cells.Item(234, 2).Value  = "Vega Monumental Concepción"
$ws.Cells.Item(234, 3).Value  = "Bíobío"
$ws.Cells.Item(234, 4).Value  = 44736
$ws.Cells.Item(234, 5).Value  = 8
$ws.Cells.Item(234, 6).Value  = "Fruta"
$ws.Cells.Item(234, 7).Value  = 100102
$ws.Cells.Item(234, 8).Value  = "Cítricos"
$ws.Cells.Item(234, 9).Value  = 100102005
$ws.Cells.Item(234, 10).Value = "Naranja"
$ws.Cells.Item(234, 11).Value = "Fukumoto"
$ws.Cells.Item(234, 12).Value = "Primera"
$ws.Cells.Item(234, 13).Value = 300
$ws.Cells.Item(234, 14).Value = 7500
$ws.Cells.Item(234, 15).Value = 8000
$ws.Cells.Item(234, 16).Value = 7750
$ws.Cells.Item(234, 17).Value = "`$/caja 15 kilos granel"
$ws.Cells.Item(234, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(234, 19).Value = 517
$ws.Cells.Item(234, 20).Value = 15

# New row 235: Thompson / Primera
$ws.Cells.Item(235, 1).Value  = 11
$ws.Cells.Item(235, 2).Value  = "Vega Monumental Concepción"
$ws.Cells.Item(235, 3).Value  = "Bíobío"
$ws.Cells.Item(235, 4).Value  = 44736
$ws.Cells.Item(235, 5).Value  = 8
$ws.Cells.Item(235, 6).Value  = "Fruta"
$ws.Cells.Item(235, 7).Value  = 100102
$ws.Cells.Item(235, 8).Value  = "Cítricos"
$ws.Cells.Item(235, 9).Value  = 100102005
$ws.Cells.Item(235, 10).Value = "Naranja"
$ws.Cells.Item(235, 11).Value = "Thompson"
$ws.Cells.Item(235, 12).Value = "Primera"
$ws.Cells.Item(235, 13).Value = 300
$ws.Cells.Item(235, 14).Value = 6500
$ws.Cells.Item(235, 15).Value = 7000
$ws.Cells.Item(235, 16).Value = 6750
$ws.Cells.Item(235, 17).Value = "`$/caja 15 kilos granel"
$ws.Cells.Item(235, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(235, 19).Value = 450
$ws.Cells.Item(235, 20).Value = 15
